# Natmi following Dr Hou advice
# Rewrites the LR-pair data block (rows 2-4 become a 2x3 ECs/sCs x ECs/FAPs/sCs
# matrix expanded to rows 2-7) with refreshed NATMI statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: ECs -> Il1b/Il1rap -> ECs ----
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il1b"
$ws.Range("C2").Value = "Il1rap"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1413.335253
$ws.Range("H2").Value = 4240.005759
$ws.Range("I2").Value = 0.9999668843963775
$ws.Range("J2").Value = 0.9999668843963775
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 20.990057
$ws.Range("N2").Value = 62.97017099999999
$ws.Range("O2").Value = 0.6163261988329278
$ws.Range("P2").Value = 0.6163261988329277
$ws.Range("Q2").Value = 29665.98752057942
$ws.Range("R2").Value = 266993.8876852148
$ws.Range("S2").Value = 0.6163057888188251
$ws.Range("T2").Value = 0.616305788818825

# ---- Row 3: ECs -> Il1b/Il1rap -> FAPs ----
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il1b"
$ws.Range("C3").Value = "Il1rap"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1413.335253
$ws.Range("H3").Value = 4240.005759
$ws.Range("I3").Value = 0.9999668843963775
$ws.Range("J3").Value = 0.9999668843963775
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.095305
$ws.Range("N3").Value = 21.285915
$ws.Range("O3").Value = 0.208337802999309
$ws.Range("P3").Value = 0.208337802999309
$ws.Range("Q3").Value = 10028.04468728717
$ws.Range("R3").Value = 90252.40218558448
$ws.Range("S3").Value = 0.2083309037672053
$ws.Range("T3").Value = 0.2083309037672053

# ---- Row 4: ECs -> Il1b/Il1rap -> sCs ----
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il1b"
$ws.Range("C4").Value = "Il1rap"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1413.335253
$ws.Range("H4").Value = 4240.005759
$ws.Range("I4").Value = 0.9999668843963775
$ws.Range("J4").Value = 0.9999668843963775
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.971371333333334
$ws.Range("N4").Value = 17.914114
$ws.Range("O4").Value = 0.1753359981677632
$ws.Range("P4").Value = 0.1753359981677632
$ws.Range("Q4").Value = 8439.549614153615
$ws.Range("R4").Value = 75955.94652738252
$ws.Range("S4").Value = 0.1753301918103471
$ws.Range("T4").Value = 0.1753301918103471

# ---- Row 5: sCs -> Il1b/Il1rap -> ECs (new row) ----
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Il1b"
$ws.Range("C5").Value = "Il1rap"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.04680500000000001
$ws.Range("H5").Value = 0.140415
$ws.Range("I5").Value = [double]"3.311560362258399E-05"
$ws.Range("J5").Value = [double]"3.311560362258399E-05"
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 20.990057
$ws.Range("N5").Value = 62.97017099999999
$ws.Range("O5").Value = 0.6163261988329278
$ws.Range("P5").Value = 0.6163261988329277
$ws.Range("Q5").Value = 0.982439617885
$ws.Range("R5").Value = 8.841956560965
$ws.Range("S5").Value = [double]"2.041001410276512E-05"
$ws.Range("T5").Value = [double]"2.041001410276512E-05"

# ---- Row 6: sCs -> Il1b/Il1rap -> FAPs (new row) ----
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Il1b"
$ws.Range("C6").Value = "Il1rap"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.04680500000000001
$ws.Range("H6").Value = 0.140415
$ws.Range("I6").Value = [double]"3.311560362258399E-05"
$ws.Range("J6").Value = [double]"3.311560362258399E-05"
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.095305
$ws.Range("N6").Value = 21.285915
$ws.Range("O6").Value = 0.208337802999309
$ws.Range("P6").Value = 0.208337802999309
$ws.Range("Q6").Value = 0.332095750525
$ws.Range("R6").Value = 2.988861754725
$ws.Range("S6").Value = [double]"6.899232103725108E-06"
$ws.Range("T6").Value = [double]"6.899232103725107E-06"

# ---- Row 7: sCs -> Il1b/Il1rap -> sCs (new row) ----
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Il1b"
$ws.Range("C7").Value = "Il1rap"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.04680500000000001
$ws.Range("H7").Value = 0.140415
$ws.Range("I7").Value = [double]"3.311560362258399E-05"
$ws.Range("J7").Value = [double]"3.311560362258399E-05"
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.971371333333334
$ws.Range("N7").Value = 17.914114
$ws.Range("O7").Value = 0.1753359981677632
$ws.Range("P7").Value = 0.1753359981677632
$ws.Range("Q7").Value = 0.2794900352566667
$ws.Range("R7").Value = 2.51541031731
$ws.Range("S7").Value = [double]"5.80635741609376E-06"
$ws.Range("T7").Value = [double]"5.80635741609376E-06"
